$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.202.67"
$ws.Range("E2").Value = "  -0.02%  "
$ws.Range("D3").Value = "1.850.20"
$ws.Range("E3").Value = "  -0.29%  "
$ws.Range("D4").Value = "'0.9997"
$ws.Range("E4").Value = "  -0.09%  "
$ws.Range("E5").Value = "  +2.01%  "
$ws.Range("D6").Value = "'0.6981"
$ws.Range("E6").Value = "  -0.58%  "
$ws.Range("E7").Value = "  -0.08%  "
$ws.Range("D8").Value = "'0.07719"
$ws.Range("E8").Value = "  +0.02%  "
$ws.Range("E9").Value = "  -0.82%  "
$ws.Range("D10").Value = "'23.50"
$ws.Range("E10").Value = "  -1.03%  "
$ws.Range("D11").Value = "'0.07822"
$ws.Range("E11").Value = "  +0.26%  "
$ws.Range("E12").Value = "  +1.19%  "
$ws.Range("D13").Value = "1.849.55"
$ws.Range("E13").Value = "  -0.50%  "
$ws.Range("D14").Value = "'5.127"
$ws.Range("E14").Value = "  +0.78%  "
$ws.Range("D15").Value = "'0.6857"
$ws.Range("E15").Value = "  +0.09%  "
$ws.Range("D16").Value = "'6.632"
$ws.Range("E16").Value = "  +2.25%  "
$ws.Range("B17").Value = "ShibaInu"
$ws.Range("C17").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D17").Value = "'0.000008309"
$ws.Range("E17").Value = "  -0.98%  "
$ws.Range("B18").Value = "WrappedBTC"
$ws.Range("C18").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D18").Value = "29.185.92"
$ws.Range("E18").Value = "  -0.11%  "
$ws.Range("D19").Value = "'241.61"
$ws.Range("E19").Value = "  -2.94%  "
$ws.Range("D20").Value = "2.087.59"
$ws.Range("E20").Value = "  -1.80%  "
$ws.Range("D21").Value = "'12.73"
$ws.Range("E21").Value = "  -0.66%  "
$ws.Range("E22").Value = "  -0.04%  "
$ws.Range("D23").Value = "'7.521"
$ws.Range("E23").Value = "  +0.14%  "
$ws.Range("E24").Value = "  -0.11%  "
$ws.Range("D25").Value = "'0.1512"
$ws.Range("E25").Value = "  -0.86%  "
$ws.Range("D26").Value = "'159.40"
$ws.Range("E26").Value = "  -0.42%  "
$ws.Range("D27").Value = "'8.827"
$ws.Range("E27").Value = "  +0.02%  "
$ws.Range("D28").Value = "'18.28"
$ws.Range("E28").Value = "  -1.10%  "
$ws.Range("D29").Value = "'1.541"
$ws.Range("E29").Value = "  -1.15%  "
$ws.Range("D30").Value = "'4.226"
$ws.Range("E30").Value = "  -0.11%  "
$ws.Range("D31").Value = "'4.179"
$ws.Range("E31").Value = "  -0.53%  "
$ws.Range("D32").Value = "'1.198"
$ws.Range("E32").Value = "  +0.24%  "
$ws.Range("D33").Value = "'0.05121"
$ws.Range("E33").Value = "  -1.28%  "
$ws.Range("D34").Value = "'0.7935"
$ws.Range("E34").Value = "  +4.87%  "
$ws.Range("E35").Value = "  +1.42%  "
$ws.Range("D36").Value = "'1.148"
$ws.Range("E36").Value = "  -0.93%  "
$ws.Range("E37").Value = "  -0.53%  "
$ws.Range("D38").Value = "1.320.51"
$ws.Range("E38").Value = "  +7.72%  "
$ws.Range("D39").Value = "'0.01874"
$ws.Range("E39").Value = "  +0.89%  "
$ws.Range("D40").Value = "'2.712"
$ws.Range("E40").Value = "  -0.45%  "
$ws.Range("D41").Value = "'0.9494"
$ws.Range("E41").Value = "  +6.04%  "
$ws.Range("D42").Value = "'6.004"
$ws.Range("E42").Value = "  +7.53%  "
$ws.Range("D43").Value = "'107.24"
$ws.Range("E43").Value = "  -2.16%  "
$ws.Range("E44").Value = "  +0.01%  "
$ws.Range("D45").Value = "'9.746"
$ws.Range("E45").Value = "  +2.73%  "
$ws.Range("D46").Value = "1.988.90"
$ws.Range("E46").Value = "  -1.17%  "
$ws.Range("D47").Value = "'0.5181"
$ws.Range("D48").Value = "'64.14"
$ws.Range("E48").Value = "  -1.52%  "
$ws.Range("D49").Value = "'1.765"
$ws.Range("E49").Value = "  +1.22%  "
$ws.Range("D50").Value = "'0.00000000119"
$ws.Range("E50").Value = "  -1.76%  "
$ws.Range("D51").Value = "'6.998"
$ws.Range("E51").Value = "  +0.23%  "
